# ELT-3A schedule fix: shift afternoon/evening slots by one 50-min turn
# so that the day covers 6 full hours per turn (commit: "6 hours by turn fix").
#
# The lunch ("Almoço") / break ("Intervalo") rows and the one-off lecture
# ("Jorge Aquino-Motores de aplicação") all move down one row, the time
# column is recomputed, and three brand-new rows are appended at the
# bottom of the table (16:40, 17:30 and a trailing 18:20 row marking the
# end of the last turn).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: the morning lecture slot reverts to the default "-" ---
$ws.Range("E4").Value = "-"

# --- Row 6 is unchanged (still the lecture slot) ---

# --- Row 7: the lecture moves here (one turn later) ---
$ws.Range("B7").Value = "Jorge Aquino-Motores de aplicação"

# --- Row 8: lunch no longer starts here ---
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"

# --- Row 9: time moves earlier (12:20) and lunch now starts here ---
$ws.Range("A9").Value = "12:20"
$ws.Range("B9").Value = "Almoço"
$ws.Range("C9").Value = "Almoço"
$ws.Range("D9").Value = "Almoço"
$ws.Range("E9").Value = "Almoço"
$ws.Range("F9").Value = "Almoço"

# --- Row 10: time shifts earlier, content ("-") unchanged ---
$ws.Range("A10").Value = "13:00"

# --- Row 11: time shifts earlier, content ("-") unchanged ---
$ws.Range("A11").Value = "13:50"

# --- Row 12: time shifts earlier, the break no longer occurs here ---
$ws.Range("A12").Value = "14:40"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"

# --- Row 13: time shifts earlier, break moves here ---
$ws.Range("A13").Value = "15:30"
$ws.Range("B13").Value = "Intervalo"
$ws.Range("C13").Value = "Intervalo"
$ws.Range("D13").Value = "Intervalo"
$ws.Range("E13").Value = "Intervalo"
$ws.Range("F13").Value = "Intervalo"

# --- Row 14: time shifts earlier; content becomes the default "-" ---
$ws.Range("A14").Value = "15:50"
$ws.Range("B14").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"

# --- Row 15 (new): 16:40 turn, default "-" across the week ---
$ws.Range("A15").Value = "16:40"
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"

# --- Row 16 (new): 17:30 turn, default "-" across the week ---
$ws.Range("A16").Value = "17:30"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"

# --- Row 17 (new): closing time marker, blank activity cells ---
$ws.Range("A17").Value = "18:20"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
